$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{Row=2;  C=4254.42;  D=128.92;              E=4898.61;  F=148.44},
    @{Row=3;  C=4704.82;  D=142.57;              E=531.65;   F=16.11},
    @{Row=4;  C=450.4;    D=13.65;               E=-4366.96; F=-132.33},
    @{Row=5;  C=4493.85;  D=136.18;              E=4996.53;  F=151.41},
    @{Row=6;  C=921.47;   D=27.92;               E=1171.22;  F=35.49},
    @{Row=7;  C=-3572.38; D=-108.26;             E=-3825.31; F=-115.92},
    @{Row=8;  C=5065.24;  D=153.49;              E=4682.85;  F=141.9},
    @{Row=9;  C=2787.75;  D=84.48;               E=2899.46;  F=87.86},
    @{Row=10; C=-2277.49; D=-69.01000000000001;  E=-1783.39; F=-54.04},
    @{Row=11; C=1373.73;  D=41.63;               E=2634.37;  F=79.83},
    @{Row=12; C=2674.61;  D=81.05;               E=1766.54;  F=53.53},
    @{Row=13; C=1300.88;  D=39.42;               E=-867.83;  F=-26.3},
    @{Row=14; C=808.88;   D=24.51;               E=887.59;   F=26.9},
    @{Row=15; C=930.39;   D=28.19;               E=740.04;   F=22.43},
    @{Row=16; C=121.51;   D=3.68;                E=-147.55;  F=-4.47}
)

foreach ($entry in $data) {
    $r = $entry.Row
    $ws.Cells.Item($r, 3).Value = $entry.C
    $ws.Cells.Item($r, 4).Value = $entry.D
    $ws.Cells.Item($r, 5).Value = $entry.E
    $ws.Cells.Item($r, 6).Value = $entry.F
}
